# issue #5: stock data from json to db
#
# The "股票" (stock) sheet gets three new metadata columns that track where
# each scraped row came from: category, source_file and index.
# "category" is inserted right after the existing "property_category"
# column (pushing date / legislator_name / legislator_id one column to the
# right), while "source_file" and "index" are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # 股票 (stock)

# Insert a new column before I; this shifts the old I/J/K (date /
# legislator_name / legislator_id) columns one to the right, to J/K/L.
$ws.Columns("I:I").Insert()

# New header + value for the inserted column.
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"

# Two brand-new columns appended after the existing data.
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"
$ws.Range("M2").Value = "tmp72b91"
$ws.Range("N2").Value = 70

# Match the look of the surrounding header/data cells by copying their
# formatting onto the newly added ones.
$ws.Range("K1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K2").Copy()
$ws.Range("M2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
